$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (Wins, Losses, Ties) in row 1, columns AD, AE, AF.
# Copy formatting from the existing header cell (AC1) so the new headers
# match the bold / bordered / centered style used by the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every player row (2-55): 81 wins, 81 losses, 0 ties.
$ws.Range("AD2:AD55").Value = 81
$ws.Range("AE2:AE55").Value = 81
$ws.Range("AF2:AF55").Value = 0
